$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 and A3 to share the new string "abc" (replacing "1-PC" / "2-PC")
$ws.Range("A2").Value = "abc"
$ws.Range("A3").Value = "abc"

# Update the selected cell/range shown in the sheet view
$ws.Range("D9").Select()
